$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("D1").Value = "ARQUIVO"
$ws.Range("E1").Value = "FORMATO"
$ws.Range("F1").Value = "DISCIPLINA"
$ws.Range("G1").Value = "TIPO DE DOCUMENTO"
$ws.Range("H1").Value = "PROPÓSITO"
$ws.Range("I1").Value = "CAMINHO DATABOOK"

# Row 2
$ws.Range("A2").Value = "documento_pid"
$ws.Range("C2").Value = "Documento PID de Teste"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "A4"
$ws.Range("F2").Value = "PROCESSO"
$ws.Range("G2").Value = "PID"
$ws.Range("H2").Value = "Para Construção"
$ws.Range("I2").Value = "DATA BOOK C&M"

# Row 3
$ws.Range("A3").Value = "documento_rir"
$ws.Range("C3").Value = "Documento RIR de Teste"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "A3"
$ws.Range("F3").Value = "INSTRUMENTAÇÃO"
$ws.Range("G3").Value = "RIR"
$ws.Range("H3").Value = "Para Construção"
$ws.Range("I3").Value = "DATA BOOK C&M"
